$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currency_movements")

# Insert a new column before column D ("currency") to hold the new "fees" column
$ws.Columns("D").Insert()

# The newly inserted column inherited formatting from its left neighbor;
# the data cell in row 2 should stay unformatted like the rest of the column
$ws.Range("D2").ClearFormats()

# Header for the new "fees" column
$ws.Range("D1").Value = "fees"

# Fees values for the four data rows
$ws.Range("D2").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0

# amounts for the last two rows are reduced by the fee amount (100 -> 98, keeping sign)
$ws.Range("C4").Value = -98
$ws.Range("C5").Value = -98
